$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.843.48'
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '1.634.75'
$ws.Range("E3").Value = '  -1.52%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5027'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.33%  '
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06408'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07695'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.55%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.636.98'
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.245'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.36%  '
$ws.Range("D14").Value = '1.859.00'
$ws.Range("E14").Value = '  -1.57%  '
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").Value = '0.0₅7924'
$ws.Range("E16").Value = '  -1.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.61'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").Value = '25.856.38'
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '203.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.323'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.942'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.978'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.929'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.32'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1146'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.71'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.706'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.241'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05012'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.82%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.271'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.185'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.49%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.535'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.353'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.83%  '
$ws.Range("D36").Value = '1.175.85'
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.8923'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.616'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5587'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01561'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.19%  '
$ws.Range("B41").Value = 'mCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.546'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.98%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.002'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.663'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8072'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '99.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.22%  '
$ws.Range("D46").Value = '1.771.36'
$ws.Range("E46").Value = '  -1.50%  '
$ws.Range("E47").Value = '  +0.04%  '
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.005'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("E50").Value = '  -2.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05036'
$ws.Range("D51").Style = "Normal"
